$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New observation rows appended at the bottom of the sheet (rows 20-25).
# Text-like columns get an explicit text NumberFormat right before the value
# is written so that numeric-looking strings (ids, dates, times, counts, ...)
# are not auto-converted by Excel into numbers or dates.

# Row 20
$ws.Range("A20").Value = 131106801
$ws.Range("B20").Value = 57300
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = 'NT'
$ws.Range("E20").Value = 102961
$ws.Range("F20").NumberFormat = "@"
$ws.Range("F20").Value = 'Drillsnäppa'
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = 'Actitis hypoleucos'
$ws.Range("H20").NumberFormat = "@"
$ws.Range("H20").Value = 'Linnaeus, 1758'
$ws.Range("I20").NumberFormat = "@"
$ws.Range("I20").Value = '2'
$ws.Range("P20").NumberFormat = "@"
$ws.Range("P20").Value = 'Paljack, Mpd'
$ws.Range("Q20").Value = 600200
$ws.Range("R20").Value = 6973062
$ws.Range("S20").Value = 10
$ws.Range("T20").NumberFormat = "@"
$ws.Range("T20").Value = 'Västernorrland'
$ws.Range("U20").NumberFormat = "@"
$ws.Range("U20").Value = 'Sundsvall'
$ws.Range("V20").NumberFormat = "@"
$ws.Range("V20").Value = 'Medelpad'
$ws.Range("W20").NumberFormat = "@"
$ws.Range("W20").Value = 'Liden'
$ws.Range("X20").NumberFormat = "@"
$ws.Range("X20").Value = '2025_0375'
$ws.Range("Y20").NumberFormat = "@"
$ws.Range("Y20").Value = '2025-06-24'
$ws.Range("Z20").NumberFormat = "@"
$ws.Range("Z20").Value = '15:05'
$ws.Range("AA20").NumberFormat = "@"
$ws.Range("AA20").Value = '2025-06-24'
$ws.Range("AB20").NumberFormat = "@"
$ws.Range("AB20").Value = '15:05'
$ws.Range("AC20").NumberFormat = "@"
$ws.Range("AC20").Value = 'par i häckbiotop'
$ws.Range("AD20").Value = $false
$ws.Range("AE20").Value = $false
$ws.Range("AG20").Value = $false
$ws.Range("AT20").NumberFormat = "@"
$ws.Range("AT20").Value = ''
$ws.Range("AW20").NumberFormat = "@"
$ws.Range("AW20").Value = 'David Isaksson'
$ws.Range("AX20").NumberFormat = "@"
$ws.Range("AX20").Value = 'Anders Forsberg'
$ws.Range("AY20").NumberFormat = "@"
$ws.Range("AY20").Value = 'Kustpaketet'

# Row 21
$ws.Range("A21").Value = 131106774
$ws.Range("B21").Value = 79833
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = 'NT'
$ws.Range("E21").Value = 229821
$ws.Range("F21").NumberFormat = "@"
$ws.Range("F21").Value = 'Vedflamlav'
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = 'Ramboldia elabens'
$ws.Range("H21").NumberFormat = "@"
$ws.Range("H21").Value = '(Fr.) Kantvilas & Elix'
$ws.Range("I21").NumberFormat = "@"
$ws.Range("I21").Value = '2'
$ws.Range("J21").NumberFormat = "@"
$ws.Range("J21").Value = 'dm²'
$ws.Range("P21").NumberFormat = "@"
$ws.Range("P21").Value = 'Nordansjö, Mpd'
$ws.Range("Q21").Value = 600324
$ws.Range("R21").Value = 6972295
$ws.Range("S21").Value = 10
$ws.Range("T21").NumberFormat = "@"
$ws.Range("T21").Value = 'Västernorrland'
$ws.Range("U21").NumberFormat = "@"
$ws.Range("U21").Value = 'Sundsvall'
$ws.Range("V21").NumberFormat = "@"
$ws.Range("V21").Value = 'Medelpad'
$ws.Range("W21").NumberFormat = "@"
$ws.Range("W21").Value = 'Liden'
$ws.Range("X21").NumberFormat = "@"
$ws.Range("X21").Value = '2025_0402'
$ws.Range("Y21").NumberFormat = "@"
$ws.Range("Y21").Value = '2025-06-25'
$ws.Range("Z21").NumberFormat = "@"
$ws.Range("Z21").Value = '08:16'
$ws.Range("AA21").NumberFormat = "@"
$ws.Range("AA21").Value = '2025-06-25'
$ws.Range("AB21").NumberFormat = "@"
$ws.Range("AB21").Value = '08:16'
$ws.Range("AD21").Value = $false
$ws.Range("AE21").Value = $false
$ws.Range("AG21").Value = $false
$ws.Range("AT21").NumberFormat = "@"
$ws.Range("AT21").Value = ''
$ws.Range("AW21").NumberFormat = "@"
$ws.Range("AW21").Value = 'David Isaksson'
$ws.Range("AX21").NumberFormat = "@"
$ws.Range("AX21").Value = 'David Isaksson'
$ws.Range("AY21").NumberFormat = "@"
$ws.Range("AY21").Value = 'Kustpaketet'

# Row 22
$ws.Range("A22").Value = 131106800
$ws.Range("B22").Value = 57300
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = 'NT'
$ws.Range("E22").Value = 102961
$ws.Range("F22").NumberFormat = "@"
$ws.Range("F22").Value = 'Drillsnäppa'
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = 'Actitis hypoleucos'
$ws.Range("H22").NumberFormat = "@"
$ws.Range("H22").Value = 'Linnaeus, 1758'
$ws.Range("I22").NumberFormat = "@"
$ws.Range("I22").Value = '2'
$ws.Range("P22").NumberFormat = "@"
$ws.Range("P22").Value = 'Paljack, Mpd'
$ws.Range("Q22").Value = 600203
$ws.Range("R22").Value = 6973054
$ws.Range("S22").Value = 10
$ws.Range("T22").NumberFormat = "@"
$ws.Range("T22").Value = 'Västernorrland'
$ws.Range("U22").NumberFormat = "@"
$ws.Range("U22").Value = 'Sundsvall'
$ws.Range("V22").NumberFormat = "@"
$ws.Range("V22").Value = 'Medelpad'
$ws.Range("W22").NumberFormat = "@"
$ws.Range("W22").Value = 'Liden'
$ws.Range("X22").NumberFormat = "@"
$ws.Range("X22").Value = '2025_0376'
$ws.Range("Y22").NumberFormat = "@"
$ws.Range("Y22").Value = '2025-06-24'
$ws.Range("Z22").NumberFormat = "@"
$ws.Range("Z22").Value = '15:06'
$ws.Range("AA22").NumberFormat = "@"
$ws.Range("AA22").Value = '2025-06-24'
$ws.Range("AB22").NumberFormat = "@"
$ws.Range("AB22").Value = '15:06'
$ws.Range("AC22").NumberFormat = "@"
$ws.Range("AC22").Value = 'par i häckbiotop'
$ws.Range("AD22").Value = $false
$ws.Range("AE22").Value = $false
$ws.Range("AG22").Value = $false
$ws.Range("AT22").NumberFormat = "@"
$ws.Range("AT22").Value = ''
$ws.Range("AW22").NumberFormat = "@"
$ws.Range("AW22").Value = 'David Isaksson'
$ws.Range("AX22").NumberFormat = "@"
$ws.Range("AX22").Value = 'Anders Forsberg'
$ws.Range("AY22").NumberFormat = "@"
$ws.Range("AY22").Value = 'Kustpaketet'

# Row 23
$ws.Range("A23").Value = 131106775
$ws.Range("B23").Value = 79862
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = 'NT'
$ws.Range("E23").Value = 6453
$ws.Range("F23").NumberFormat = "@"
$ws.Range("F23").Value = 'Vedskivlav'
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = 'Hertelidea botryosa'
$ws.Range("H23").NumberFormat = "@"
$ws.Range("H23").Value = '(Fr.) Printzen & Kantvilas'
$ws.Range("I23").NumberFormat = "@"
$ws.Range("I23").Value = '5'
$ws.Range("J23").NumberFormat = "@"
$ws.Range("J23").Value = 'cm²'
$ws.Range("P23").NumberFormat = "@"
$ws.Range("P23").Value = 'Nordansjö, Mpd'
$ws.Range("Q23").Value = 600380
$ws.Range("R23").Value = 6972329
$ws.Range("S23").Value = 10
$ws.Range("T23").NumberFormat = "@"
$ws.Range("T23").Value = 'Västernorrland'
$ws.Range("U23").NumberFormat = "@"
$ws.Range("U23").Value = 'Sundsvall'
$ws.Range("V23").NumberFormat = "@"
$ws.Range("V23").Value = 'Medelpad'
$ws.Range("W23").NumberFormat = "@"
$ws.Range("W23").Value = 'Liden'
$ws.Range("X23").NumberFormat = "@"
$ws.Range("X23").Value = '2025_0401'
$ws.Range("Y23").NumberFormat = "@"
$ws.Range("Y23").Value = '2025-06-25'
$ws.Range("Z23").NumberFormat = "@"
$ws.Range("Z23").Value = '08:07'
$ws.Range("AA23").NumberFormat = "@"
$ws.Range("AA23").Value = '2025-06-25'
$ws.Range("AB23").NumberFormat = "@"
$ws.Range("AB23").Value = '08:07'
$ws.Range("AD23").Value = $false
$ws.Range("AE23").Value = $false
$ws.Range("AG23").Value = $false
$ws.Range("AT23").NumberFormat = "@"
$ws.Range("AT23").Value = ''
$ws.Range("AW23").NumberFormat = "@"
$ws.Range("AW23").Value = 'David Isaksson'
$ws.Range("AX23").NumberFormat = "@"
$ws.Range("AX23").Value = 'David Isaksson, Karin Halldin'
$ws.Range("AY23").NumberFormat = "@"
$ws.Range("AY23").Value = 'Kustpaketet'

# Row 24
$ws.Range("A24").Value = 131106670
$ws.Range("B24").Value = 79833
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = 'NT'
$ws.Range("E24").Value = 229821
$ws.Range("F24").NumberFormat = "@"
$ws.Range("F24").Value = 'Vedflamlav'
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = 'Ramboldia elabens'
$ws.Range("H24").NumberFormat = "@"
$ws.Range("H24").Value = '(Fr.) Kantvilas & Elix'
$ws.Range("I24").NumberFormat = "@"
$ws.Range("I24").Value = '4'
$ws.Range("J24").NumberFormat = "@"
$ws.Range("J24").Value = 'dm²'
$ws.Range("P24").NumberFormat = "@"
$ws.Range("P24").Value = 'Nordansjö, Mpd'
$ws.Range("Q24").Value = 600304
$ws.Range("R24").Value = 6972313
$ws.Range("S24").Value = 10
$ws.Range("T24").NumberFormat = "@"
$ws.Range("T24").Value = 'Västernorrland'
$ws.Range("U24").NumberFormat = "@"
$ws.Range("U24").Value = 'Sundsvall'
$ws.Range("V24").NumberFormat = "@"
$ws.Range("V24").Value = 'Medelpad'
$ws.Range("W24").NumberFormat = "@"
$ws.Range("W24").Value = 'Liden'
$ws.Range("X24").NumberFormat = "@"
$ws.Range("X24").Value = '2025_0506'
$ws.Range("Y24").NumberFormat = "@"
$ws.Range("Y24").Value = '2025-06-25'
$ws.Range("Z24").NumberFormat = "@"
$ws.Range("Z24").Value = '08:28'
$ws.Range("AA24").NumberFormat = "@"
$ws.Range("AA24").Value = '2025-06-25'
$ws.Range("AB24").NumberFormat = "@"
$ws.Range("AB24").Value = '08:28'
$ws.Range("AD24").Value = $false
$ws.Range("AE24").Value = $false
$ws.Range("AG24").Value = $false
$ws.Range("AT24").NumberFormat = "@"
$ws.Range("AT24").Value = ''
$ws.Range("AW24").NumberFormat = "@"
$ws.Range("AW24").Value = 'David Isaksson'
$ws.Range("AX24").NumberFormat = "@"
$ws.Range("AX24").Value = 'Karin Halldin'
$ws.Range("AY24").NumberFormat = "@"
$ws.Range("AY24").Value = 'Kustpaketet'

# Row 25
$ws.Range("A25").Value = 131106669
$ws.Range("B25").Value = 57073
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = 'LC'
$ws.Range("E25").Value = 100138
$ws.Range("F25").NumberFormat = "@"
$ws.Range("F25").Value = 'Tjäder'
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = 'Tetrao urogallus'
$ws.Range("H25").NumberFormat = "@"
$ws.Range("H25").Value = 'Linnaeus, 1758'
$ws.Range("I25").NumberFormat = "@"
$ws.Range("I25").Value = '1'
$ws.Range("P25").NumberFormat = "@"
$ws.Range("P25").Value = 'Nordansjö, Mpd'
$ws.Range("Q25").Value = 600375
$ws.Range("R25").Value = 6972262
$ws.Range("S25").Value = 10
$ws.Range("T25").NumberFormat = "@"
$ws.Range("T25").Value = 'Västernorrland'
$ws.Range("U25").NumberFormat = "@"
$ws.Range("U25").Value = 'Sundsvall'
$ws.Range("V25").NumberFormat = "@"
$ws.Range("V25").Value = 'Medelpad'
$ws.Range("W25").NumberFormat = "@"
$ws.Range("W25").Value = 'Liden'
$ws.Range("X25").NumberFormat = "@"
$ws.Range("X25").Value = '2025_0507'
$ws.Range("Y25").NumberFormat = "@"
$ws.Range("Y25").Value = '2025-06-25'
$ws.Range("Z25").NumberFormat = "@"
$ws.Range("Z25").Value = '08:41'
$ws.Range("AA25").NumberFormat = "@"
$ws.Range("AA25").Value = '2025-06-25'
$ws.Range("AB25").NumberFormat = "@"
$ws.Range("AB25").Value = '08:41'
$ws.Range("AD25").Value = $false
$ws.Range("AE25").Value = $false
$ws.Range("AG25").Value = $false
$ws.Range("AT25").NumberFormat = "@"
$ws.Range("AT25").Value = ''
$ws.Range("AW25").NumberFormat = "@"
$ws.Range("AW25").Value = 'David Isaksson'
$ws.Range("AX25").NumberFormat = "@"
$ws.Range("AX25").Value = 'Karin Halldin'
$ws.Range("AY25").NumberFormat = "@"
$ws.Range("AY25").Value = 'Kustpaketet'
